$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Placement, Rack) after "Office / Department" (col M),
# before the old "Unit Cost" column (old col N). This shifts Unit Cost, Total
# Cost, Physical Condition, Remarks two columns to the right.
$ws.Columns("N:O").Insert()

# New header cells for the inserted columns.
$ws.Range("N2").Value = "Placement"
$ws.Range("O2").Value = "Rack"

# --- Row 3 updates ---
$ws.Range("B3").Value = "FFE-TRA-BCD-1011"
$ws.Range("D3").Value = "testing"
$ws.Range("E3").Value = 121
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 1
$ws.Range("N3").Value = "test place"
$ws.Range("O3").Value = "test rack"
$ws.Range("P3").Value = ""
$ws.Range("R3").Value = "12 USD"

# --- Row 4 updates ---
$ws.Range("B4").Value = "FFE-TRA-BCD-1013"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "sss"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 2
$ws.Range("J4").Value = ""
$ws.Range("N4").Value = "test place"
$ws.Range("O4").Value = "test rack"
$ws.Range("P4").Value = ""
$ws.Range("Q4").Value = "0.00 "
$ws.Range("R4").Value = "0 "

# --- Protected ranges now need to cover the two extra columns (A:S instead
# of A:Q) for each data row. ---
$ws.Protection.AllowEditRanges.Add("p6f5710c8199129451d13d7cdd4fbfe8e", $ws.Range("A3:S3"), "C724")
$ws.Protection.AllowEditRanges.Add("p5789a9ca630561449388a885f4dbc751", $ws.Range("A4:S4"), "C724")
